# Update imputed values in the RandomForest result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.0879
$ws.Range("B8").Value = 6.493099999999999
$ws.Range("B10").Value = 5.4922
$ws.Range("B12").Value = 6.297499999999999
$ws.Range("C12").Value = -11.8838
$ws.Range("C15").Value = -13.4777
$ws.Range("C17").Value = -13.5966
$ws.Range("B18").Value = 6.828999999999997
$ws.Range("C26").Value = -13.24550000000001
$ws.Range("C27").Value = -12.6251
$ws.Range("C28").Value = -13.3785
$ws.Range("B37").Value = 8.767900000000003
$ws.Range("C37").Value = -12.9425
$ws.Range("C47").Value = -12.54489999999999
$ws.Range("B55").Value = 6.246299999999998
$ws.Range("C65").Value = -12.8563
$ws.Range("B68").Value = 4.759999999999994
$ws.Range("C73").Value = -11.33180000000001
$ws.Range("B77").Value = 9.371400000000005
$ws.Range("B78").Value = 9.810799999999995
$ws.Range("B81").Value = 5.702400000000003
$ws.Range("B82").Value = 5.548499999999999
$ws.Range("C84").Value = -12.977
$ws.Range("C85").Value = -12.4267
$ws.Range("C93").Value = -10.2696
$ws.Range("C95").Value = -12.1191
$ws.Range("C98").Value = -12.426
$ws.Range("C99").Value = -11.96200000000001
$ws.Range("C101").Value = -13.3491
